# Update workbook for "Add data for 2022-07-16"
# - Advances the reporting cutoff from 2022-07-07 to 2022-07-08
# - Updates the sheet name and the "current month" header label accordingly
# - Updates the carjacking counts for every "July" column (one per year in the
#   rolling dataset) to reflect the extra day (July 8) of historical data now
#   included in each year's July total.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab
$ws.Name = "Through 2022-07-08"

# Update the shared header label for the current (in-progress) month column
$ws.Range("B1").Value = "July 2022 (through July 08)"

# Cell updates: Range -> new value
$updates = @{
    "I2"   = 5
    "AD2"  = 4
    "AD3"  = 1
    "P5"   = 2
    "B6"   = 4
    "AR6"  = 2
    "P8"   = 5
    "AY16" = 1
    "I18"  = 2
    "AD18" = 2
    "I22"  = 1
    "W26"  = 1
    "B27"  = 2
    "B29"  = 2
    "I36"  = 1
    "AK47" = 1
    "AR49" = 1
    "I50"  = 1
    "AK52" = 2
    "P53"  = 3
    "AD62" = 1
    "AD65" = 1
    "I96"  = 4
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
